$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column S by copying the formatting (and values, to be
# overwritten below) from column R, rows 3-8, which carries the visual
# style (borders / number format) that the new 2022 column should use.
$ws.Range("R3:R8").Copy($ws.Range("S3:S8"))

# --- New column S (year 2022) ---
$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 13.600365850576139
$ws.Range("S6").Value = 9.2742414863791556
$ws.Range("S7").Value = 17.303523954725925
$ws.Range("S8").Value = 205.5

# --- Updated figures for existing years (recalculated source data) ---
$ws.Range("P5").Value = 23.111083656771282
$ws.Range("Q5").Value = 24.08077930418019
$ws.Range("R5").Value = 19.336931533747723

$ws.Range("P6").Value = 14.322631450320875
$ws.Range("Q6").Value = 13.073459110725862
$ws.Range("R6").Value = 10.464141365743002

$ws.Range("P7").Value = 23.612622725489956

# --- Selection moves to Q15 (matches the sheetView selection in the diff) ---
$ws.Range("Q15").Select()
